$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1112.1305
$ws.Range("I15").Value = 1112.1305
$ws.Range("K15").Value = 3336.3915
$ws.Range("M15").Value = -3167.3915

$ws.Range("H17").Value = 1470.7755
$ws.Range("J17").Value = 1470.7755
$ws.Range("L17").Value = 4412.3265
$ws.Range("N17").Value = -4748.3265

$ws.Range("H62").Value = 948.75
$ws.Range("I62").Value = 897.5
$ws.Range("K62").Value = 897.5
$ws.Range("M62").Value = -273.5

$ws.Range("H65").Value = 948.75
$ws.Range("I65").Value = 897.5
$ws.Range("K65").Value = 4487.5
$ws.Range("M65").Value = -1367.5

$ws.Range("H74").Value = 5523.385
$ws.Range("I74").Value = 5433.222
$ws.Range("K74").Value = 5433.222
$ws.Range("M74").Value = -4497.222

$ws.Range("H77").Value = 5523.385
$ws.Range("I77").Value = 5433.222
$ws.Range("K77").Value = 27166.11
$ws.Range("M77").Value = -22486.11

$ws.Range("H138").Value = 17372.2
$ws.Range("I138").Value = 53941.684
$ws.Range("J138").Value = 2267.413
$ws.Range("K138").Value = 161825.052
$ws.Range("L138").Value = 6802.239
$ws.Range("M138").Value = -156685.052
$ws.Range("N138").Value = -17082.239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5954325
$ws.Range("I45").Value = 2252.3
$ws.Range("J45").Value = 11365300
$ws.Range("K45").Value = 2252.3
$ws.Range("L45").Value = 11365300
$ws.Range("M45").Value = -1875.3
$ws.Range("N45").Value = -11366054

$ws.Range("H61").Value = 42562.16
$ws.Range("J61").Value = 127817.625
$ws.Range("L61").Value = 127817.625
$ws.Range("N61").Value = -128241.625

$ws.Range("H134").Value = 71547.5
$ws.Range("J134").Value = 71547.5
$ws.Range("L134").Value = 71547.5
$ws.Range("N134").Value = -81687.5

$ws.Range("H136").Value = 42562.16
$ws.Range("J136").Value = 127817.625
$ws.Range("L136").Value = 383452.875
$ws.Range("N136").Value = -388552.875

$ws.Range("H138").Value = 68562.39999999999
$ws.Range("J138").Value = 67105.5
$ws.Range("L138").Value = 67105.5
$ws.Range("N138").Value = -77385.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 380978.7
$ws.Range("I99").Value = 63838.875
$ws.Range("J99").Value = 771304.6
$ws.Range("K99").Value = 63838.875
$ws.Range("L99").Value = 771304.6
$ws.Range("M99").Value = -62340.875
$ws.Range("N99").Value = -774300.6

$ws.Range("H134").Value = 4515
$ws.Range("J134").Value = 7670
$ws.Range("L134").Value = 23010
$ws.Range("N134").Value = -28080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 16628.4
$ws.Range("J28").Value = 16628.4
$ws.Range("L28").Value = 16628.4
$ws.Range("N28").Value = -17118.4

$ws.Range("H86").Value = 7942.75
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 7942.75
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H132").Value = 1705.4138
$ws.Range("I132").Value = 1461.2727
$ws.Range("K132").Value = 4383.8181
$ws.Range("M132").Value = -1853.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 799.7083
$ws.Range("I5").Value = 702
$ws.Range("J5").Value = 915.1818
$ws.Range("K5").Value = 2106
$ws.Range("L5").Value = 2745.5454
$ws.Range("M5").Value = -1994
$ws.Range("N5").Value = -2969.5454

$ws.Range("H34").Value = 3694.1667
$ws.Range("I34").Value = 585
$ws.Range("J34").Value = 5248.75
$ws.Range("K34").Value = 1755
$ws.Range("L34").Value = 15746.25
$ws.Range("M34").Value = -1671
$ws.Range("N34").Value = -15914.25

$ws.Range("H128").Value = 424237.25
$ws.Range("I128").Value = 424237.25
$ws.Range("K128").Value = 1272711.75
$ws.Range("M128").Value = -1267731.75

$ws.Range("H135").Value = 799.7083
$ws.Range("I135").Value = 702
$ws.Range("J135").Value = 915.1818
$ws.Range("K135").Value = 6318
$ws.Range("L135").Value = 8236.636199999999
$ws.Range("M135").Value = -3783
$ws.Range("N135").Value = -13306.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 97585.60000000001
$ws.Range("I70").Value = 49797.727
$ws.Range("K70").Value = 49797.727
$ws.Range("M70").Value = -49527.727

$ws.Range("H73").Value = 97585.60000000001
$ws.Range("I73").Value = 49797.727
$ws.Range("K73").Value = 49797.727
$ws.Range("M73").Value = -48861.727

$ws.Range("H135").Value = 53804.094
$ws.Range("J135").Value = 53804.094
$ws.Range("L135").Value = 53804.094
$ws.Range("N135").Value = -63944.094

$ws.Range("H140").Value = 98557.875
$ws.Range("J140").Value = 98557.875
$ws.Range("L140").Value = 98557.875
$ws.Range("N140").Value = -108917.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4084.1614
$ws.Range("J22").Value = 8117.643
$ws.Range("L22").Value = 8117.643
$ws.Range("N22").Value = -8707.643

$ws.Range("H27").Value = 4084.1614
$ws.Range("J27").Value = 8117.643
$ws.Range("L27").Value = 8117.643
$ws.Range("N27").Value = -8331.643

$ws.Range("H46").Value = 6390.2334
$ws.Range("I46").Value = 11614.9
$ws.Range("K46").Value = 11614.9
$ws.Range("M46").Value = -11426.9

$ws.Range("H132").Value = 4329.9
$ws.Range("J132").Value = 4412.5
$ws.Range("L132").Value = 13237.5
$ws.Range("N132").Value = -18297.5

$ws.Range("H134").Value = 130289.6
$ws.Range("J134").Value = 130289.6
$ws.Range("L134").Value = 130289.6
$ws.Range("N134").Value = -140429.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1786.1111
$ws.Range("I107").Value = 775.25
$ws.Range("J107").Value = 2074.9285
$ws.Range("K107").Value = 2325.75
$ws.Range("L107").Value = 6224.7855
$ws.Range("M107").Value = -405.75
$ws.Range("N107").Value = -10064.7855

$ws.Range("H137").Value = 131999.5
$ws.Range("J137").Value = 131999.5
$ws.Range("L137").Value = 131999.5
$ws.Range("N137").Value = -142199.5

$ws.Range("H140").Value = 149831.5
$ws.Range("J140").Value = 149831.5
$ws.Range("L140").Value = 149831.5
$ws.Range("N140").Value = -160191.5

$ws.Range("H141").Value = 58750.4
$ws.Range("J141").Value = 58750.4
$ws.Range("L141").Value = 58750.4
$ws.Range("N141").Value = -69110.39999999999
